# Adds the K-Means data to the "Sheet 1" worksheet (three new rows),
# matching the commit "adds the K-Means data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Insert three new data rows right after the header row (row 2), pushing
# the existing KNN / PAM Regress rows down.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(5).Insert()

# --- Copy cell formatting from the existing table rows so the new rows
# visually match (band-header style for the first new row, continuation
# style for the other two). This also avoids the default "Text" number
# format that Insert() applies to brand-new rows.
$ws.Range("A6:F6").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)

$ws.Range("A7:F7").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$ws.Range("A7:F7").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Values are written in the same order the original author entered them
# (Kmeans/K=22, KMeans/K=28, Kmeans/K=30) before sorting the table by
# algorithm name, so that new shared-string entries come out in the same
# order as the authoritative workbook.

# --- Row 3: Kmeans, K=22 ---------------------------------------------
$ws.Range("A3").Value = "Kmeans"
$ws.Range("B3").Value = "forestfires.csv"
$ws.Range("C3").Value = "K=22"
$ws.Range("D3").Value = 4176.542345
$ws.Range("E3").Formula = "=SQRT(D3)"
$ws.Range("F3").Formula = "=E3/63.655818"
$ws.Range("F3").NumberFormat = "General"

# --- Row 5: KMeans, K=28 ----------------------------------------------
$ws.Range("A5").Value = "KMeans"
$ws.Range("B5").Value = "forestfires.csv"
$ws.Range("C5").Value = "K=28"
$ws.Range("D5").Value = 4178.456212
$ws.Range("E5").Formula = "=SQRT(D5)"
$ws.Range("F5").Formula = "=E5/63.655818"
$ws.Range("F5").NumberFormat = "General"

# --- Row 4: Kmeans, K=30 ----------------------------------------------
$ws.Range("A4").Value = "Kmeans"
$ws.Range("B4").Value = "forestfires.csv"
$ws.Range("C4").Value = "K=30"
$ws.Range("D4").Value = 4174.986352
$ws.Range("E4").Formula = "=SQRT(D4)"
$ws.Range("F4").Formula = "=E4/63.655818"
$ws.Range("F4").NumberFormat = "General"
